$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(18).Delete()
for ($r = 1; $r -le 22; $r++) {
    $t = $ws.Cells.Item($r, 1).Text
    $n = $ws.Cells.Item($r, 2).Text
    Write-Host "$r : $t | $n"
}
